# first sweep cleaning data columns to conform to specs--done by chase
#
# 1) Normalize the two category labels so they use underscores instead of
#    spaces (matches the shared-string edits in the diff).
# 2) Update the sheet's view/selection state: scroll so column E is the
#    left-most visible column and select N12 (matches the sheetView edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Data column cleanup -------------------------------------------------
$ws.Cells.Replace("Environmental Perturbation", "Environmental_Perturbation")
$ws.Cells.Replace("KN99 alpha", "KN99_alpha")

# --- 2. View / selection state ----------------------------------------------
$ws.Range("N12").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
